$d = $word.ActiveDocument

$d.Content.Find.Execute("740÷7=105, 5", $true, $false, $false, $false, $false, $true, 1, $false, "293÷5=58, 3", 2) | Out-Null
$d.Content.Find.Execute("110÷4=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "472÷6=78, 4", 2) | Out-Null
$d.Content.Find.Execute("319÷4=79, 3", $true, $false, $false, $false, $false, $true, 1, $false, "400÷3=133, 1", 2) | Out-Null
$d.Content.Find.Execute("268÷7=38, 2", $true, $false, $false, $false, $false, $true, 1, $false, "343÷5=68, 3", 2) | Out-Null
$d.Content.Find.Execute("725÷3=241, 2", $true, $false, $false, $false, $false, $true, 1, $false, "807÷7=115, 2", 2) | Out-Null
$d.Content.Find.Execute("866÷3=288, 2", $true, $false, $false, $false, $false, $true, 1, $false, "571÷6=95, 1", 2) | Out-Null
$d.Content.Find.Execute("625÷7=89, 2", $true, $false, $false, $false, $false, $true, 1, $false, "685÷6=114, 1", 2) | Out-Null
$d.Content.Find.Execute("869÷9=96, 5", $true, $false, $false, $false, $false, $true, 1, $false, "335÷9=37, 2", 2) | Out-Null
$d.Content.Find.Execute("534÷4=133, 2", $true, $false, $false, $false, $false, $true, 1, $false, "423÷6=70, 3", 2) | Out-Null
$d.Content.Find.Execute("697÷6=116, 1", $true, $false, $false, $false, $false, $true, 1, $false, "341÷4=85, 1", 2) | Out-Null
$d.Content.Find.Execute("196÷3=65, 1", $true, $false, $false, $false, $false, $true, 1, $false, "100÷3=33, 1", 2) | Out-Null
$d.Content.Find.Execute("793÷2=396, 1", $true, $false, $false, $false, $false, $true, 1, $false, "282÷5=56, 2", 2) | Out-Null
$d.Content.Find.Execute("214÷5=42, 4", $true, $false, $false, $false, $false, $true, 1, $false, "347÷6=57, 5", 2) | Out-Null
$d.Content.Find.Execute("840÷3=280, 0", $true, $false, $false, $false, $false, $true, 1, $false, "573÷4=143, 1", 2) | Out-Null
$d.Content.Find.Execute("760÷8=95, 0", $true, $false, $false, $false, $false, $true, 1, $false, "967÷5=193, 2", 2) | Out-Null
$d.Content.Find.Execute("125÷7=17, 6", $true, $false, $false, $false, $false, $true, 1, $false, "586÷5=117, 1", 2) | Out-Null
$d.Content.Find.Execute("699÷4=174, 3", $true, $false, $false, $false, $false, $true, 1, $false, "980÷2=490, 0", 2) | Out-Null
$d.Content.Find.Execute("874÷4=218, 2", $true, $false, $false, $false, $false, $true, 1, $false, "456÷9=50, 6", 2) | Out-Null
$d.Content.Find.Execute("431÷9=47, 8", $true, $false, $false, $false, $false, $true, 1, $false, "549÷7=78, 3", 2) | Out-Null
$d.Content.Find.Execute("468÷8=58, 4", $true, $false, $false, $false, $false, $true, 1, $false, "169÷6=28, 1", 2) | Out-Null
$d.Content.Find.Execute("455÷3=151, 2", $true, $false, $false, $false, $false, $true, 1, $false, "373÷6=62, 1", 2) | Out-Null
$d.Content.Find.Execute("143÷3=47, 2", $true, $false, $false, $false, $false, $true, 1, $false, "531÷7=75, 6", 2) | Out-Null
$d.Content.Find.Execute("291÷6=48, 3", $true, $false, $false, $false, $false, $true, 1, $false, "688÷2=344, 0", 2) | Out-Null
$d.Content.Find.Execute("526÷6=87, 4", $true, $false, $false, $false, $false, $true, 1, $false, "621÷5=124, 1", 2) | Out-Null
$d.Content.Find.Execute("791÷8=98, 7", $true, $false, $false, $false, $false, $true, 1, $false, "224÷9=24, 8", 2) | Out-Null
